$d = $word.ActiveDocument

# Locate the Use Case Scenarios bullet that describes the confirmation page
# and mark it (and a new trailing note) in red, the same way a reviewer
# would select the bullet and type the extra "NOT APPLICABLE" caveat.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*The confirmation page includes a thank-you message*") {

        $pStart = $p.Range.Start
        $textLen = $p.Range.Text.Length          # includes the trailing paragraph mark
        $sentenceEnd = $pStart + $textLen - 1     # position right before the paragraph mark

        $newText = " -> NOT APPLICABLE FOR THIS WEBSITE."

        # Insert the new sentence right after the existing one, before the
        # paragraph mark, as its own run.
        $insPoint = $d.Range($sentenceEnd, $sentenceEnd)
        $insPoint.InsertAfter($newText)
        $newEnd = $sentenceEnd + $newText.Length

        # Color just the newly-inserted text first so it ends up as its own
        # run, separate from the original sentence.
        $rNew = $d.Range($sentenceEnd, $newEnd)
        $rNew.Font.Color = 255

        # Now color the whole paragraph (original sentence + new sentence +
        # paragraph mark). This both re-colors the original sentence and
        # stamps the paragraph mark's run properties (w:pPr/w:rPr), matching
        # how Word records "whole paragraph" character formatting.
        $whole = $p.Range
        $whole.Font.Color = 255

        break
    }
}
